# ---------------------------------------------------------------------------
# Applies the "minor improvements to two PowerPoint slides" commit:
#   * Slide 13 ("Common Conventions for Rules")
#       - Swap/expand the "Set T ..." / "Set N ..." bullet paragraphs.
#       - Move the bordered "Using these conventions ..." box down.
#   * Slide 14 ("Example: Grammar for CPRL")
#       - Drop the optional "-" sign from the constDecl EBNF rule.
#   * Slide 39 ("Associativity")
#       - Remove the word "adjacent" from the first bullet.
#       - Delete the separate "Note: All operators in CPRL ..." text box.
#   * Slide 5 (CPRL code fragment)
#       - "proc main" -> "proc main()"
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 13: rewrite the "Set T" / "Set N" bullets (reordered + extended) and
# move the framed callout box down to make room for the extra line of text.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$body13 = $s13.Shapes.Item(4)          # "Rectangle 3" - bulleted body text
$tr13 = $body13.TextFrame.TextRange

$oldSetT = "Set T consists of all terminal symbols appearing in the rules."
$oldSetN = "Set N consists of all nonterminals appearing in the rules."

# The 2nd paragraph (currently the "Set T ..." sentence) becomes the new
# "Set N ..." sentence. "nonterminals" is marked off with a placeholder
# ("@@") so a later pass can retouch just that word, which is how
# PowerPoint ends up splitting it into its own run.
$fullText = $tr13.Text
$idxSetT = $fullText.IndexOf($oldSetT)
$rangeP2 = $tr13.Characters($idxSetT + 1, $oldSetT.Length)
$rangeP2.Text = "Set N consists of all @@ appearing in the rules; i.e., any symbol that appears on the left side of a rule."

# The 3rd paragraph (currently the "Set N ..." sentence) becomes the new
# "Set T ..." sentence.
$fullText = $tr13.Text
$idxSetN = $fullText.IndexOf($oldSetN)
$rangeP3 = $tr13.Characters($idxSetN + 1, $oldSetN.Length)
$rangeP3.Text = "Set T consists of all terminal symbols appearing in the rules; i.e., any symbol that does not appear on the left side of a rule."

# Retouch the placeholder so "nonterminals" becomes its own run.
$fullText = $tr13.Text
$phIdx = $fullText.IndexOf("@@")
$placeholder = $tr13.Characters($phIdx + 1, 2)
$placeholder.Text = "nonterminals"

# Move the framed "Using these conventions ..." box further down the slide.
$callout13 = $s13.Shapes.Item(5)        # "Rectangle 4"
$callout13.Top = 378                    # 4800600 EMU

# ---------------------------------------------------------------------------
# Slide 14: simplify the constDecl grammar rule (drop the optional sign).
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$body14 = $s14.Shapes.Item(2)           # "Rectangle 3" - EBNF grammar text
$tr14 = $body14.TextFrame.TextRange

$full14 = $tr14.Text
$oldRule = ' ":=" [ "-" ] literal ";" .'
$newRule = ' ":=" literal ";" .'
$idxRule = $full14.IndexOf($oldRule)
$ruleRange = $tr14.Characters($idxRule + 1, $oldRule.Length)
$ruleRange.Text = $newRule

# ---------------------------------------------------------------------------
# Slide 39: trim the "Associativity" explanation and drop the redundant note.
# ---------------------------------------------------------------------------
$s39 = $p.Slides.Item(39)
$body39 = $s39.Shapes.Item(2)           # "Content Placeholder 2"
$tr39 = $body39.TextFrame.TextRange

$full39 = $tr39.Text
$oldSentence = "Specifies the evaluation order of adjacent operators with the same precedence when there are no parentheses."
$newSentence = "Specifies the evaluation order of operators with the same precedence when there are no parentheses."
$idx39 = $full39.IndexOf($oldSentence)
$sentenceRange = $tr39.Characters($idx39 + 1, $oldSentence.Length)
$sentenceRange.Text = $newSentence

# Delete the standalone "Note: All operators in CPRL are left associative." box.
$noteBox = $s39.Shapes.Item(5)          # "TextBox 5"
$noteBox.Delete()

# ---------------------------------------------------------------------------
# Slide 5: complete the procedure signature with an empty parameter list.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$code5 = $s5.Shapes.Item(6)             # "Rectangle 12"
$tr5 = $code5.TextFrame.TextRange

$full5 = $tr5.Text
$oldProc = "proc main"
$newProc = "proc main()"
$idx5 = $full5.IndexOf($oldProc)
$procRange = $tr5.Characters($idx5 + 1, $oldProc.Length)
$procRange.Text = $newProc
